# Auto-generated Excel COM-interop script
# Applies numeric cell updates across multiple worksheets
# as described by the source XML diff (commit: chore: update Sheets via scheduled runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5685544.5
$ws.Range("I51").Value = 22731210
$ws.Range("J51").Value = 3656
$ws.Range("K51").Value = 22731210
$ws.Range("L51").Value = 3656
$ws.Range("M51").Value = -22730726
$ws.Range("N51").Value = -4624
$ws.Range("H135").Value = 992.44446
$ws.Range("I135").Value = 548.1053000000001
$ws.Range("K135").Value = 4932.947700000001
$ws.Range("M135").Value = -2397.947700000001
$ws.Range("H138").Value = 1568.6897
$ws.Range("I138").Value = 1249.6451
$ws.Range("J138").Value = 1935
$ws.Range("K138").Value = 3748.9353
$ws.Range("L138").Value = 5805
$ws.Range("M138").Value = 1391.0647
$ws.Range("N138").Value = -16085
$ws.Range("H141").Value = 9036
$ws.Range("I141").Value = 2900.5334
$ws.Range("J141").Value = 55052
$ws.Range("K141").Value = 8701.600199999999
$ws.Range("L141").Value = 165156
$ws.Range("M141").Value = -3521.600199999999
$ws.Range("N141").Value = -175516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 44995
$ws.Range("J55").Value = 44995
$ws.Range("L55").Value = 44995
$ws.Range("N55").Value = -45625
$ws.Range("H109").Value = 24000
$ws.Range("J109").Value = 24000
$ws.Range("L109").Value = 24000
$ws.Range("N109").Value = -26774
$ws.Range("H112").Value = 16972.4
$ws.Range("J112").Value = 16972.4
$ws.Range("L112").Value = 16972.4
$ws.Range("N112").Value = -19926.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2368.65
$ws.Range("I16").Value = 2564.0667
$ws.Range("J16").Value = 1782.4
$ws.Range("K16").Value = 2564.0667
$ws.Range("L16").Value = 1782.4
$ws.Range("M16").Value = -2277.0667
$ws.Range("N16").Value = -2356.4
$ws.Range("H22").Value = 554.0526
$ws.Range("I22").Value = 541.6667
$ws.Range("J22").Value = 600.5
$ws.Range("K22").Value = 541.6667
$ws.Range("L22").Value = 600.5
$ws.Range("M22").Value = -191.6667
$ws.Range("N22").Value = -1300.5
$ws.Range("H31").Value = 2640.2144
$ws.Range("I31").Value = 2010.8667
$ws.Range("J31").Value = 3366.3845
$ws.Range("K31").Value = 2010.8667
$ws.Range("L31").Value = 3366.3845
$ws.Range("M31").Value = -1715.8667
$ws.Range("N31").Value = -3956.3845
$ws.Range("H34").Value = 2640.2144
$ws.Range("I34").Value = 2010.8667
$ws.Range("J34").Value = 3366.3845
$ws.Range("K34").Value = 2010.8667
$ws.Range("L34").Value = 3366.3845
$ws.Range("M34").Value = -1808.8667
$ws.Range("N34").Value = -3770.3845
$ws.Range("H58").Value = 8668.799999999999
$ws.Range("I58").Value = 4314.75
$ws.Range("J58").Value = 13644.857
$ws.Range("K58").Value = 4314.75
$ws.Range("L58").Value = 13644.857
$ws.Range("M58").Value = -4111.75
$ws.Range("N58").Value = -14050.857
$ws.Range("H99").Value = 2320.889
$ws.Range("I99").Value = 1605.1428
$ws.Range("J99").Value = 2776.3635
$ws.Range("K99").Value = 1605.1428
$ws.Range("L99").Value = 2776.3635
$ws.Range("M99").Value = -107.1428000000001
$ws.Range("N99").Value = -5772.363499999999
$ws.Range("H113").Value = 2368.65
$ws.Range("I113").Value = 2564.0667
$ws.Range("J113").Value = 1782.4
$ws.Range("K113").Value = 2564.0667
$ws.Range("L113").Value = 1782.4
$ws.Range("M113").Value = -394.0666999999999
$ws.Range("N113").Value = -6122.4
$ws.Range("H126").Value = 2320.889
$ws.Range("I126").Value = 1605.1428
$ws.Range("J126").Value = 2776.3635
$ws.Range("K126").Value = 4815.428400000001
$ws.Range("L126").Value = 8329.0905
$ws.Range("M126").Value = -2345.428400000001
$ws.Range("N126").Value = -13269.0905
$ws.Range("H132").Value = 2768.1765
$ws.Range("I132").Value = 2137.2632
$ws.Range("J132").Value = 3567.3333
$ws.Range("K132").Value = 6411.7896
$ws.Range("L132").Value = 10701.9999
$ws.Range("M132").Value = -3881.7896
$ws.Range("N132").Value = -15761.9999
$ws.Range("H134").Value = 2552.1538
$ws.Range("I134").Value = 1822.25
$ws.Range("J134").Value = 2876.5557
$ws.Range("K134").Value = 5466.75
$ws.Range("L134").Value = 8629.667099999999
$ws.Range("M134").Value = -2931.75
$ws.Range("N134").Value = -13699.6671
$ws.Range("H136").Value = 8668.799999999999
$ws.Range("I136").Value = 4314.75
$ws.Range("J136").Value = 13644.857
$ws.Range("K136").Value = 12944.25
$ws.Range("L136").Value = 40934.571
$ws.Range("M136").Value = -10394.25
$ws.Range("N136").Value = -46034.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5393
$ws.Range("I64").Value = 786
$ws.Range("K64").Value = 2358
$ws.Range("M64").Value = -2088
$ws.Range("H67").Value = 5393
$ws.Range("I67").Value = 786
$ws.Range("K67").Value = 2358
$ws.Range("M67").Value = -1422
$ws.Range("H131").Value = 997.7954999999999
$ws.Range("I131").Value = 499.84616
$ws.Range("J131").Value = 1206.6129
$ws.Range("K131").Value = 1499.53848
$ws.Range("L131").Value = 3619.8387
$ws.Range("M131").Value = 3540.46152
$ws.Range("N131").Value = -13699.8387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2253.84
$ws.Range("I102").Value = 1605.8572
$ws.Range("J102").Value = 3078.5454
$ws.Range("K102").Value = 1605.8572
$ws.Range("L102").Value = 3078.5454
$ws.Range("M102").Value = 16.14280000000008
$ws.Range("N102").Value = -6322.5454
$ws.Range("H122").Value = 3022.25
$ws.Range("I122").Value = 2525.7273
$ws.Range("J122").Value = 3343.5293
$ws.Range("K122").Value = 7577.1819
$ws.Range("L122").Value = 10030.5879
$ws.Range("M122").Value = -5127.1819
$ws.Range("N122").Value = -14930.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1250760
$ws.Range("J46").Value = 2500897.5
$ws.Range("L46").Value = 2500897.5
$ws.Range("N46").Value = -2501273.5
$ws.Range("H55").Value = 364.5
$ws.Range("I55").Value = 334.35715
$ws.Range("J55").Value = 470
$ws.Range("K55").Value = 334.35715
$ws.Range("L55").Value = 470
$ws.Range("M55").Value = -161.35715
$ws.Range("N55").Value = -816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2298.6875
$ws.Range("I122").Value = 1461.091
$ws.Range("J122").Value = 4141.4
$ws.Range("K122").Value = 4383.272999999999
$ws.Range("L122").Value = 12424.2
$ws.Range("M122").Value = -1933.272999999999
$ws.Range("N122").Value = -17324.2
$ws.Range("H132").Value = 1908.5106
$ws.Range("I132").Value = 1612.96
$ws.Range("K132").Value = 4838.88
$ws.Range("M132").Value = -2308.88

Write-Host "Updated 165 cells across 7 worksheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)."
